# ---------------------------------------------------------------------------
# Updates "广州-漫展信息.xlsx" (Guangzhou convention-info workbook) to match
# the newer scrape: a new exhibition ("广东·广州EY动漫嘉年华") is inserted,
# several "想去人数" (want-to-go) counters are bumped, and one performance's
# counter is bumped too.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $text) {
    # Force text storage (column holds values that look numeric/date-like,
    # e.g. "60" or "2024.02.12", which Excel would otherwise auto-coerce).
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value2 = $text
    $c.Style = "Normal"
}

function Insert-EventRow($ws, $atRow, $idx, $date, $name, $place, $timerange, $want, $price, $stage, $link) {
    $ws.Rows.Item($atRow).Insert()

    # Match column-A's bold/bordered numbering style used by every other row.
    $ws.Cells.Item($atRow - 1, 1).Copy()
    $ws.Cells.Item($atRow, 1).PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($atRow, 1).Value2 = $idx
    Set-TextCell $ws $atRow 2 $date
    $ws.Cells.Item($atRow, 3).Value2 = $name
    $ws.Cells.Item($atRow, 4).Value2 = $place
    $ws.Cells.Item($atRow, 5).Value2 = $timerange
    $ws.Cells.Item($atRow, 6).Value2 = $want
    Set-TextCell $ws $atRow 7 $price
    $ws.Cells.Item($atRow, 8).Value2 = $stage
    $ws.Cells.Item($atRow, 9).Value2 = $link
}

# ===========================================================================
# Sheet "展览" (exhibitions)
# ===========================================================================
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Cells.Item(2, 6).Value2 = 55
$ws1.Cells.Item(3, 6).Value2 = 978
$ws1.Cells.Item(5, 6).Value2 = 10871
$ws1.Cells.Item(6, 6).Value2 = 1153
$ws1.Cells.Item(7, 6).Value2 = 361
$ws1.Cells.Item(8, 6).Value2 = 594
$ws1.Cells.Item(9, 6).Value2 = 1945
$ws1.Cells.Item(10, 6).Value2 = 588
$ws1.Cells.Item(11, 6).Value2 = 732
$ws1.Cells.Item(12, 6).Value2 = 227
$ws1.Cells.Item(13, 6).Value2 = 306
$ws1.Cells.Item(14, 6).Value2 = 273
$ws1.Cells.Item(15, 6).Value2 = 275
$ws1.Cells.Item(16, 6).Value2 = 988
$ws1.Cells.Item(17, 6).Value2 = 377

Insert-EventRow $ws1 18 17 "2024.02.12" "广东·广州EY动漫嘉年华" `
    "机场路1399号广州百信广场二期 李宁运动中心" "2024.02.12 10:00-02.12 17:00" `
    2 "63" $false `
    "https://show.bilibili.com/platform/detail.html?id=80574&msource=Msearch_colligation"

# Existing rows 18-24 shifted down to 19-25; the A-column index sequence
# (0,1,2,...) continues unbroken, and "want to go" counts are bumped too.
$ws1.Cells.Item(19, 1).Value2 = 18
$ws1.Cells.Item(19, 6).Value2 = 201
$ws1.Cells.Item(20, 1).Value2 = 19
$ws1.Cells.Item(20, 6).Value2 = 421
$ws1.Cells.Item(21, 1).Value2 = 20
$ws1.Cells.Item(21, 6).Value2 = 641
$ws1.Cells.Item(22, 1).Value2 = 21
$ws1.Cells.Item(22, 6).Value2 = 784
$ws1.Cells.Item(23, 1).Value2 = 22
$ws1.Cells.Item(23, 6).Value2 = 177
$ws1.Cells.Item(24, 1).Value2 = 23
$ws1.Cells.Item(24, 6).Value2 = 448
$ws1.Cells.Item(25, 1).Value2 = 24
$ws1.Cells.Item(25, 6).Value2 = 186

# ===========================================================================
# Sheet "演出" (performances)
# ===========================================================================
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(7, 6).Value2 = 630

# ===========================================================================
# Sheet "全部类型" (all types combined)
# ===========================================================================
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Cells.Item(4, 6).Value2 = 55
$ws4.Cells.Item(5, 6).Value2 = 978
$ws4.Cells.Item(8, 6).Value2 = 10871
$ws4.Cells.Item(9, 6).Value2 = 1153
$ws4.Cells.Item(11, 6).Value2 = 361
$ws4.Cells.Item(12, 6).Value2 = 594
$ws4.Cells.Item(13, 6).Value2 = 1945
$ws4.Cells.Item(14, 6).Value2 = 588
$ws4.Cells.Item(15, 6).Value2 = 732
$ws4.Cells.Item(17, 6).Value2 = 227
$ws4.Cells.Item(18, 6).Value2 = 306
$ws4.Cells.Item(19, 6).Value2 = 273
$ws4.Cells.Item(20, 6).Value2 = 275
$ws4.Cells.Item(21, 6).Value2 = 988
$ws4.Cells.Item(22, 6).Value2 = 377
$ws4.Cells.Item(23, 6).Value2 = 630

Insert-EventRow $ws4 24 23 "2024.02.12" "广东·广州EY动漫嘉年华" `
    "机场路1399号广州百信广场二期 李宁运动中心" "2024.02.12 10:00-02.12 17:00" `
    2 "63" $false `
    "https://show.bilibili.com/platform/detail.html?id=80574&msource=Msearch_colligation"

# Existing rows 24-31 shifted down to 25-32; the A-column index sequence
# (0,1,2,...) continues unbroken, and "want to go" counts are bumped too.
$ws4.Cells.Item(25, 1).Value2 = 24
$ws4.Cells.Item(25, 6).Value2 = 201
$ws4.Cells.Item(26, 1).Value2 = 25
$ws4.Cells.Item(26, 6).Value2 = 421
$ws4.Cells.Item(27, 1).Value2 = 26
$ws4.Cells.Item(27, 6).Value2 = 641
$ws4.Cells.Item(28, 1).Value2 = 27
$ws4.Cells.Item(28, 6).Value2 = 784
$ws4.Cells.Item(29, 1).Value2 = 28
$ws4.Cells.Item(30, 1).Value2 = 29
$ws4.Cells.Item(30, 6).Value2 = 177
$ws4.Cells.Item(31, 1).Value2 = 30
$ws4.Cells.Item(31, 6).Value2 = 448
$ws4.Cells.Item(32, 1).Value2 = 31
$ws4.Cells.Item(32, 6).Value2 = 186

Write-Host "Edit complete."
